$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.050.87'
$ws.Range("E2").Value = '  -2.34%  '
$ws.Range("D3").Value = '1.818.22'
$ws.Range("E3").Value = '  -1.54%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.20%  '
$ws.Range("D5").Value = "'310.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.94%  '
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("D7").Value = "'0.4216"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.13%  '
$ws.Range("D8").Value = "'0.3667"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.87%  '
$ws.Range("D9").Value = "'0.07209"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.95%  '
$ws.Range("D10").Value = "'0.8428"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.90%  '
$ws.Range("D11").Value = "'20.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.63%  '
$ws.Range("D12").Value = '1.834.56'
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").Value = "'6.634"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").Value = "'0.07075"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.76%  '
$ws.Range("D15").Value = "'5.268"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.19%  '
$ws.Range("D16").Value = "'88.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.26%  '
$ws.Range("D18").Value = "'0.000008811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("D20").Value = "'14.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.14%  '
$ws.Range("D21").Value = '27.108.38'
$ws.Range("D22").Value = "'5.113"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.08%  '
$ws.Range("D23").Value = "'10.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("D24").Value = '2.042.63'
$ws.Range("E24").Value = '  -2.12%  '
$ws.Range("D25").Value = "'1.979"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("D26").Value = "'151.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("D27").Value = "'2.240"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.23%  '
$ws.Range("D28").Value = "'18.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.87%  '
$ws.Range("D29").Value = "'5.197"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.44%  '
$ws.Range("D30").Value = "'115.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.38%  '
$ws.Range("D31").Value = "'0.08796"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("D32").Value = "'1.178"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.19%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = "'2.975"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.68%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = "'0.7384"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.85%  '
$ws.Range("D35").Value = "'4.410"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.14%  '
$ws.Range("D36").Value = "'1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").Value = "'1.094"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.71%  '
$ws.Range("D38").Value = "'0.01963"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("D39").Value = "'0.05247"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("D40").Value = "'7.257"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("D41").Value = "'2.872"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("D42").Value = "'0.1688"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.21%  '
$ws.Range("D43").Value = "'0.5023"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.18%  '
$ws.Range("D44").Value = "'8.551"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.57%  '
$ws.Range("D45").Value = "'10.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.72%  '
$ws.Range("D46").Value = "'0.4748"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").Value = "'106.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.25%  '
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("D49").Value = "'0.06370"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.92%  '
$ws.Range("D50").Value = "'1.649"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.57%  '
$ws.Range("D51").Value = "'1.877"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.57%  '
